$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (K2:T2)
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.09207700000000001
$ws.Range("N2").Value = 0.184154
$ws.Range("O2").Value = 0.0789959771480734
$ws.Range("P2").Value = 0.05545240531440215
$ws.Range("Q2").Value = 0.006512207209666668
$ws.Range("R2").Value = 0.039073243258
$ws.Range("S2").Value = 0.0789959771480734
$ws.Range("T2").Value = 0.05545240531440215

# Update row 3 values (O3, P3, S3, T3)
$ws.Range("O3").Value = 0.8491451975864605
$ws.Range("P3").Value = 0.8941052196698643
$ws.Range("S3").Value = 0.8491451975864605
$ws.Range("T3").Value = 0.8941052196698643

# Update row 4 values (M4:T4)
$ws.Range("M4").Value = 0.083758
$ws.Range("N4").Value = 0.167516
$ws.Range("O4").Value = 0.07185882526546619
$ws.Range("P4").Value = 0.05044237501573352
$ws.Range("Q4").Value = 0.005923840388666667
$ws.Range("R4").Value = 0.035543042332
$ws.Range("S4").Value = 0.07185882526546619
$ws.Range("T4").Value = 0.05044237501573352

# Delete row 5 entirely (the last data row, referencing "Neutrophils")
$ws.Rows("5:5").Delete()
